$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.960.97"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "1.741.06"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'249.94"
$ws.Range("E5").Value = "  +6.72%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'0.5141"
$ws.Range("E7").Value = "  -2.86%  "
$ws.Range("D8").Value = "'0.2749"
$ws.Range("E8").Value = "  -1.23%  "
$ws.Range("D9").Value = "'0.06184"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "1.744.96"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").Value = "'0.07223"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "'15.11"
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("D13").Value = "'0.6492"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").Value = "'4.632"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").Value = "'77.54"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "25.997.84"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").Value = "'11.84"
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("D20").Value = "'0.000006797"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").Value = "1.966.55"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "'4.266"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").Value = "'8.672"
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("D24").Value = "'5.367"
$ws.Range("E24").Value = "  +3.02%  "
$ws.Range("D25").Value = "'136.01"
$ws.Range("E25").Value = "  -2.38%  "
$ws.Range("D26").Value = "'1.511"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").Value = "'1.781"
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("D29").Value = "'105.95"
$ws.Range("E29").Value = "  +1.71%  "
$ws.Range("D30").Value = "'3.952"
$ws.Range("E30").Value = "  +4.09%  "
$ws.Range("D31").Value = "'0.08212"
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").Value = "'3.645"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").Value = "'0.04706"
$ws.Range("E33").Value = "  +3.64%  "
$ws.Range("D34").Value = "'2.659"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("D35").Value = "'0.9951"
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("D36").Value = "'0.6232"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("D37").Value = "'2.731"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("D38").Value = "'0.01615"
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("D39").Value = "'1.912"
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("D41").Value = "'100.18"
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("D42").Value = "'0.7595"
$ws.Range("E42").Value = "  +2.72%  "
$ws.Range("D43").Value = "'0.3852"
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("D44").Value = "'5.026"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").Value = "'6.306"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").Value = "'0.1131"
$ws.Range("E46").Value = "  -1.11%  "
$ws.Range("D47").Value = "'55.60"
$ws.Range("E47").Value = "  +2.98%  "
$ws.Range("E48").Value = "  -2.19%  "
$ws.Range("D49").Value = "'30.76"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").Value = "'7.536"
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("D51").Value = "'0.3426"
$ws.Range("E51").Value = "  -1.15%  "
